$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The content of row 5 and row 6 (the species-observation records) is being
# swapped. Columns that are already identical between the two rows (C, I, P,
# Q, R, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY) need no change.

# --- Row 5 becomes the "Buskskvätta" (bird) record ---
$ws.Range("A5").Value = 111669462
$ws.Range("B5").Value = 56890
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 102995
$ws.Range("F5").Value = "Buskskvätta"
$ws.Range("G5").Value = "Saxicola rubetra"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("S5").Value = 50
$ws.Range("AC5").ClearContents()
$ws.Range("AI5").ClearContents()

# --- Row 6 becomes the "Hällebräcka" (plant) record ---
$ws.Range("A6").Value = 111669452
$ws.Range("B6").Value = 99136
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 1449
$ws.Range("F6").Value = "Hällebräcka"
$ws.Range("G6").Value = "Saxifraga osloënsis"
$ws.Range("H6").Value = "Knaben"
$ws.Range("S6").Value = 10
$ws.Range("AC6").Value = "över 1000 ex i fin blom"
$ws.Range("AI6").Value = "I kanten av markväg på kalhygge"
